$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A17:C17")

$r.Clear()  # clear contents AND formats entirely, wiping fgColor "memory"

$ws.Cells.Item(17,1).Value2 = "Tuolumne County Transit Agency (TCTA)"
$ws.Cells.Item(17,3).Value2 = "9R02-91057"

$r.Font.Size = 12
$r.Interior.Pattern = -4142  # xlPatternNone -- now with no prior fgColor memory

Write-Host "done"
